# SGI-CP.xlsx : "Actualizacion avance del cronograma del proyecto"
#
# 1) In sheet "Cronograma #1": the verification tasks for the UI-spec and
#    DB-spec documents (currently rows 51-52) are reordered to sit right
#    after "Especificar Requerimiento 8..." (now rows 48-49), pushing the
#    "Gestion de Roles / Seguimiento de Incidencias / Asignacion de
#    Personal" tasks down to rows 50-52.
# 2) The progress column (H) for rows 46-52 goes from 0% to 100%.
# 3) In sheet "Control de Versiones" a new version-log row is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cronograma #1")

function Copy-RowPartial($srcRow, $dstRow) {
    # Copy everything except column I (a merged block we must not disturb)
    $ws.Range("A$srcRow`:H$srcRow").Copy($ws.Range("A$dstRow")) | Out-Null
    $ws.Range("J$srcRow`:K$srcRow").Copy($ws.Range("J$dstRow")) | Out-Null
}

function Copy-RowToScratch($srcRow, $scratchRow) {
    $ws.Range("A$srcRow`:H$srcRow").Copy($ws.Range("M$scratchRow")) | Out-Null
    $ws.Range("J$srcRow`:K$srcRow").Copy($ws.Range("U$scratchRow")) | Out-Null
}

function Copy-ScratchToRow($scratchRow, $dstRow) {
    $ws.Range("M$scratchRow`:T$scratchRow").Copy($ws.Range("A$dstRow")) | Out-Null
    $ws.Range("U$scratchRow`:V$scratchRow").Copy($ws.Range("J$dstRow")) | Out-Null
}

# Stash the two rows that are moving up (UI-spec, BD-spec) in a scratch
# area far to the right (columns M..V) that has no data/formatting.
Copy-RowToScratch 51 900
Copy-RowToScratch 52 901

# Shift "Gestion de Roles", "Seguimiento de Incidencias" and "Asignacion
# de Personal" down by two rows (process bottom-up so we never read a row
# we already overwrote).
Copy-RowPartial 50 52
Copy-RowPartial 49 51
Copy-RowPartial 48 50

# Drop the stashed UI-spec/BD-spec rows into their new home.
Copy-ScratchToRow 900 48
Copy-ScratchToRow 901 49

# Clean the scratch area back to empty.
$ws.Range("M900:V901").Clear() | Out-Null

# All of these tasks are now finished -> 100% progress.
$ws.Range("H46:H52").Value2 = 1

# --- Control de Versiones: log the update ---
$vws = $wb.Worksheets.Item("Control de Versiones")

# Seed row 7 with row 6's formatting, then fix column A to follow the
# numeric-date style used by rows 4/5 instead of row 6's "2.0" text style.
$vws.Range("A6:E6").Copy($vws.Range("A7")) | Out-Null
$vws.Range("A4").Copy($vws.Range("A7")) | Out-Null

$vws.Range("A7").Value2 = 44929
$vws.Range("B7").Value2 = 45096
$vws.Range("C7").Value2 = "Luis Balarezo"
$vws.Range("D7").Value2 = "Actualización del avance de los ítems al 100%"
$vws.Range("E7").Value2 = $vws.Range("E6").Value2
